$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 0.522104
$ws.Range("H2").Value2 = 1.566312
$ws.Range("I2").Value2 = 0.1803483207170308
$ws.Range("J2").Value2 = 0.1803483207170308
$ws.Range("M2").Value2 = 8.554479333333333
$ws.Range("N2").Value2 = 25.663438
$ws.Range("O2").Value2 = 0.1655051910559175
$ws.Range("P2").Value2 = 0.1655051910559175
$ws.Range("Q2").Value2 = 4.466327877850667
$ws.Range("R2").Value2 = 40.196950900656
$ws.Range("S2").Value2 = 0.02984858327688607
$ws.Range("T2").Value2 = 0.02984858327688607
$ws.Range("G3").Value2 = 0.522104
$ws.Range("H3").Value2 = 1.566312
$ws.Range("I3").Value2 = 0.1803483207170308
$ws.Range("J3").Value2 = 0.1803483207170308
$ws.Range("M3").Value2 = 20.28486166666667
$ws.Range("N3").Value2 = 60.854585
$ws.Range("O3").Value2 = 0.392455200938143
$ws.Range("P3").Value2 = 0.392455200938143
$ws.Range("Q3").Value2 = 10.59080741561333
$ws.Range("R3").Value2 = 95.31726674052
$ws.Range("S3").Value2 = 0.07077863644585897
$ws.Range("T3").Value2 = 0.07077863644585897
$ws.Range("G4").Value2 = 0.522104
$ws.Range("H4").Value2 = 1.566312
$ws.Range("I4").Value2 = 0.1803483207170308
$ws.Range("J4").Value2 = 0.1803483207170308
$ws.Range("M4").Value2 = 5.037112666666666
$ws.Range("N4").Value2 = 15.111338
$ws.Range("O4").Value2 = 0.09745400763531942
$ws.Range("P4").Value2 = 0.09745400763531943
$ws.Range("Q4").Value2 = 2.629896671717333
$ws.Range("R4").Value2 = 23.669070045456
$ws.Range("S4").Value2 = 0.01757566662417456
$ws.Range("T4").Value2 = 0.01757566662417456
$ws.Range("G5").Value2 = 0.522104
$ws.Range("H5").Value2 = 1.566312
$ws.Range("I5").Value2 = 0.1803483207170308
$ws.Range("J5").Value2 = 0.1803483207170308
$ws.Range("M5").Value2 = 17.810622
$ws.Range("N5").Value2 = 53.431866
$ws.Range("O5").Value2 = 0.34458560037062
$ws.Range("P5").Value2 = 0.34458560037062
$ws.Range("Q5").Value2 = 9.298996988688
$ws.Range("R5").Value2 = 83.690972898192
$ws.Range("S5").Value2 = 0.06214543437011118
$ws.Range("T5").Value2 = 0.06214543437011118
$ws.Range("I6").Value2 = 0.4274461949996817
$ws.Range("J6").Value2 = 0.4274461949996816
$ws.Range("M6").Value2 = 8.554479333333333
$ws.Range("N6").Value2 = 25.663438
$ws.Range("O6").Value2 = 0.1655051910559175
$ws.Range("P6").Value2 = 0.1655051910559175
$ws.Range("Q6").Value2 = 10.58570908460911
$ws.Range("R6").Value2 = 95.271381761482
$ws.Range("S6").Value2 = 0.0707445641695473
$ws.Range("T6").Value2 = 0.07074456416954729
$ws.Range("I7").Value2 = 0.4274461949996817
$ws.Range("J7").Value2 = 0.4274461949996816
$ws.Range("M7").Value2 = 20.28486166666667
$ws.Range("N7").Value2 = 60.854585
$ws.Range("O7").Value2 = 0.392455200938143
$ws.Range("P7").Value2 = 0.392455200938143
$ws.Range("Q7").Value2 = 25.10142769159056
$ws.Range("S7").Value2 = 0.1677534823488447
$ws.Range("T7").Value2 = 0.1677534823488447
$ws.Range("I8").Value2 = 0.4274461949996817
$ws.Range("J8").Value2 = 0.4274461949996816
$ws.Range("M8").Value2 = 5.037112666666666
$ws.Range("N8").Value2 = 15.111338
$ws.Range("O8").Value2 = 0.09745400763531942
$ws.Range("P8").Value2 = 0.09745400763531943
$ws.Range("Q8").Value2 = 6.233156599953555
$ws.Range("R8").Value2 = 56.098409399582
$ws.Range("S8").Value2 = 0.04165634475118721
$ws.Range("T8").Value2 = 0.04165634475118721
$ws.Range("I9").Value2 = 0.4274461949996817
$ws.Range("J9").Value2 = 0.4274461949996816
$ws.Range("M9").Value2 = 17.810622
$ws.Range("N9").Value2 = 53.431866
$ws.Range("O9").Value2 = 0.34458560037062
$ws.Range("P9").Value2 = 0.34458560037062
$ws.Range("Q9").Value2 = 22.039688888286
$ws.Range("R9").Value2 = 198.357199994574
$ws.Range("S9").Value2 = 0.1472918037301024
$ws.Range("T9").Value2 = 0.1472918037301024
$ws.Range("E10").Value2 = 2
$ws.Range("F10").Value2 = 0.6666666666666666
$ws.Range("G10").Value2 = 1.135425333333333
$ws.Range("H10").Value2 = 3.406276
$ws.Range("I10").Value2 = 0.3922054842832876
$ws.Range("J10").Value2 = 0.3922054842832876
$ws.Range("M10").Value2 = 8.554479333333333
$ws.Range("N10").Value2 = 25.663438
$ws.Range("O10").Value2 = 0.1655051910559175
$ws.Range("P10").Value2 = 0.1655051910559175
$ws.Range("Q10").Value2 = 9.71297254854311
$ws.Range("R10").Value2 = 87.416752936888
$ws.Range("S10").Value2 = 0.06491204360948417
$ws.Range("T10").Value2 = 0.06491204360948419
$ws.Range("E11").Value2 = 2
$ws.Range("F11").Value2 = 0.6666666666666666
$ws.Range("G11").Value2 = 1.135425333333333
$ws.Range("H11").Value2 = 3.406276
$ws.Range("I11").Value2 = 0.3922054842832876
$ws.Range("J11").Value2 = 0.3922054842832876
$ws.Range("M11").Value2 = 20.28486166666667
$ws.Range("N11").Value2 = 60.854585
$ws.Range("O11").Value2 = 0.392455200938143
$ws.Range("P11").Value2 = 0.392455200938143
$ws.Range("Q11").Value2 = 23.03194581949555
$ws.Range("R11").Value2 = 207.28751237546
$ws.Range("S11").Value2 = 0.1539230821434393
$ws.Range("T11").Value2 = 0.1539230821434393
$ws.Range("E12").Value2 = 2
$ws.Range("F12").Value2 = 0.6666666666666666
$ws.Range("G12").Value2 = 1.135425333333333
$ws.Range("H12").Value2 = 3.406276
$ws.Range("I12").Value2 = 0.3922054842832876
$ws.Range("J12").Value2 = 0.3922054842832876
$ws.Range("M12").Value2 = 5.037112666666666
$ws.Range("N12").Value2 = 15.111338
$ws.Range("O12").Value2 = 0.09745400763531942
$ws.Range("P12").Value2 = 0.09745400763531943
$ws.Range("Q12").Value2 = 5.719265328587555
$ws.Range("R12").Value2 = 51.473387957288
$ws.Range("S12").Value2 = 0.03822199625995766
$ws.Range("T12").Value2 = 0.03822199625995767
$ws.Range("E13").Value2 = 2
$ws.Range("F13").Value2 = 0.6666666666666666
$ws.Range("G13").Value2 = 1.135425333333333
$ws.Range("H13").Value2 = 3.406276
$ws.Range("I13").Value2 = 0.3922054842832876
$ws.Range("J13").Value2 = 0.3922054842832876
$ws.Range("M13").Value2 = 17.810622
$ws.Range("N13").Value2 = 53.431866
$ws.Range("O13").Value2 = 0.34458560037062
$ws.Range("P13").Value2 = 0.34458560037062
$ws.Range("Q13").Value2 = 20.222631421224
$ws.Range("R13").Value2 = 182.003682791016
$ws.Range("S13").Value2 = 0.1351483622704064
$ws.Range("T13").Value2 = 0.1351483622704064
